# Auto-generated edit script: updates FFXIV leve profit calculation
# values (currentAveragePrice / LevePrice / LeveProfit columns) across
# all crafting class sheets to reflect refreshed market board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3092.6667
$ws.Range("J40").Value = 2214
$ws.Range("L40").Value = 2214
$ws.Range("N40").Value = -2564
$ws.Range("H52").Value = 500
$ws.Range("J52").Value = 500
$ws.Range("L52").Value = 1500
$ws.Range("N52").Value = -1820
$ws.Range("H54").Value = 3333.3333
$ws.Range("I54").Value = 3333.3333
$ws.Range("K54").Value = 3333.3333
$ws.Range("M54").Value = -2847.3333
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H88").Value = 2499.5
$ws.Range("I88").Value = 999
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 999
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -593
$ws.Range("N88").Value = -4812
$ws.Range("H91").Value = 2499.5
$ws.Range("I91").Value = 999
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 999
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = 405
$ws.Range("N91").Value = -6808
$ws.Range("H100").Value = 451.6
$ws.Range("I100").Value = 413
$ws.Range("K100").Value = 413
$ws.Range("M100").Value = 128
$ws.Range("H138").Value = 2429.5
$ws.Range("I138").Value = 1957.56
$ws.Range("K138").Value = 5872.68
$ws.Range("M138").Value = -732.6800000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4904.6177
$ws.Range("I32").Value = 2472.5
$ws.Range("K32").Value = 2472.5
$ws.Range("M32").Value = -2185.5
$ws.Range("H33").Value = 6881.5
$ws.Range("I33").Value = 2026
$ws.Range("J33").Value = 8500
$ws.Range("K33").Value = 2026
$ws.Range("L33").Value = 8500
$ws.Range("M33").Value = -1697
$ws.Range("N33").Value = -9158
$ws.Range("H56").Value = 28750
$ws.Range("J56").Value = 28750
$ws.Range("L56").Value = 28750
$ws.Range("N56").Value = -30234
$ws.Range("H61").Value = 2290.8
$ws.Range("I61").Value = 2167.8333
$ws.Range("K61").Value = 2167.8333
$ws.Range("M61").Value = -1955.8333
$ws.Range("H63").Value = 2320
$ws.Range("I63").Value = 2400
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2400
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1714
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 2320
$ws.Range("I66").Value = 2400
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 12000
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -8568
$ws.Range("N66").Value = -16864
$ws.Range("H74").Value = 689.7143
$ws.Range("I74").Value = 504.42105
$ws.Range("K74").Value = 504.42105
$ws.Range("M74").Value = 369.57895
$ws.Range("H77").Value = 689.7143
$ws.Range("I77").Value = 504.42105
$ws.Range("K77").Value = 2522.10525
$ws.Range("M77").Value = 1845.89475
$ws.Range("H97").Value = 403.07144
$ws.Range("I97").Value = 366.36
$ws.Range("K97").Value = 366.36
$ws.Range("M97").Value = 129.64
$ws.Range("H110").Value = 6946854.5
$ws.Range("I110").Value = 7939076.5
$ws.Range("K110").Value = 7939076.5
$ws.Range("M110").Value = -7937031.5
$ws.Range("H122").Value = 13332.294
$ws.Range("I122").Value = 13729.25
$ws.Range("J122").Value = 12379.6
$ws.Range("K122").Value = 41187.75
$ws.Range("L122").Value = 37138.8
$ws.Range("M122").Value = -38737.75
$ws.Range("N122").Value = -42038.8
$ws.Range("H136").Value = 2290.8
$ws.Range("I136").Value = 2167.8333
$ws.Range("K136").Value = 6503.499899999999
$ws.Range("M136").Value = -3953.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 26224.75
$ws.Range("J35").Value = 32449.5
$ws.Range("L35").Value = 32449.5
$ws.Range("N35").Value = -33069.5
$ws.Range("H82").Value = 99998.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 99998.5
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 99998.5
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -100764.5
$ws.Range("H85").Value = 99998.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 99998.5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 99998.5
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -102650.5
$ws.Range("H107").Value = 1278.1
$ws.Range("I107").Value = 1203.5555
$ws.Range("J107").Value = 1949
$ws.Range("K107").Value = 1203.5555
$ws.Range("L107").Value = 1949
$ws.Range("M107").Value = 716.4445000000001
$ws.Range("N107").Value = -5789

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 34655
$ws.Range("I3").Value = 34000
$ws.Range("J3").Value = 34982.5
$ws.Range("K3").Value = 34000
$ws.Range("L3").Value = 34982.5
$ws.Range("M3").Value = -33887
$ws.Range("N3").Value = -35208.5
$ws.Range("H69").Value = 26400
$ws.Range("J69").Value = 39500
$ws.Range("L69").Value = 39500
$ws.Range("N69").Value = -40998
$ws.Range("H72").Value = 26400
$ws.Range("J72").Value = 39500
$ws.Range("L72").Value = 118500
$ws.Range("N72").Value = -125988

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1113169.5
$ws.Range("I4").Value = 1978277.8
$ws.Range("J4").Value = 887.2857
$ws.Range("K4").Value = 5934833.4
$ws.Range("L4").Value = 2661.8571
$ws.Range("M4").Value = -5934721.4
$ws.Range("N4").Value = -2885.8571
$ws.Range("H5").Value = 1470
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H23").Value = 250068.25
$ws.Range("J23").Value = 250068.25
$ws.Range("L23").Value = 750204.75
$ws.Range("N23").Value = -750674.75
$ws.Range("H60").Value = 309.8
$ws.Range("I60").Value = 137.25
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 411.75
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -160.75
$ws.Range("N60").Value = -3502
$ws.Range("H109").Value = 1195.5
$ws.Range("I109").Value = 1195.5
$ws.Range("K109").Value = 3586.5
$ws.Range("M109").Value = -2546.5
$ws.Range("H113").Value = 605.5
$ws.Range("I113").Value = 549.5
$ws.Range("J113").Value = 633.5
$ws.Range("K113").Value = 1648.5
$ws.Range("L113").Value = 1900.5
$ws.Range("M113").Value = 521.5
$ws.Range("N113").Value = -6240.5
$ws.Range("H121").Value = 772.5
$ws.Range("I121").Value = 772.5
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 2317.5
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -1007.5
$ws.Range("N121").ClearContents()
$ws.Range("H135").Value = 1470
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3734.3333
$ws.Range("I80").Value = 2350
$ws.Range("J80").Value = 6503
$ws.Range("K80").Value = 2350
$ws.Range("L80").Value = 6503
$ws.Range("M80").Value = -1352
$ws.Range("N80").Value = -8499
$ws.Range("H83").Value = 3734.3333
$ws.Range("I83").Value = 2350
$ws.Range("J83").Value = 6503
$ws.Range("K83").Value = 11750
$ws.Range("L83").Value = 32515
$ws.Range("M83").Value = -6758
$ws.Range("N83").Value = -42499
$ws.Range("H122").Value = 60597.65
$ws.Range("I122").Value = 1440.1428
$ws.Range("J122").Value = 336666
$ws.Range("K122").Value = 4320.428400000001
$ws.Range("L122").Value = 1009998
$ws.Range("M122").Value = -1870.428400000001
$ws.Range("N122").Value = -1014898
$ws.Range("H126").Value = 5207.25
$ws.Range("J126").Value = 5566.6665
$ws.Range("L126").Value = 16699.9995
$ws.Range("N126").Value = -21639.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 171749.67
$ws.Range("I46").Value = 5499.6665
$ws.Range("J46").Value = 337999.66
$ws.Range("K46").Value = 5499.6665
$ws.Range("L46").Value = 337999.66
$ws.Range("M46").Value = -5311.6665
$ws.Range("N46").Value = -338375.66

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5324.8335
$ws.Range("I81").Value = 6487.5
$ws.Range("J81").Value = 2999.5
$ws.Range("K81").Value = 12975
$ws.Range("L81").Value = 5999
$ws.Range("M81").Value = -11914
$ws.Range("N81").Value = -8121
$ws.Range("H84").Value = 5324.8335
$ws.Range("I84").Value = 6487.5
$ws.Range("J84").Value = 2999.5
$ws.Range("K84").Value = 64875
$ws.Range("L84").Value = 29995
$ws.Range("M84").Value = -59571
$ws.Range("N84").Value = -40603
$ws.Range("H96").Value = 5999
$ws.Range("I96").Value = 5999
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 5999
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -4626
$ws.Range("N96").ClearContents()
$ws.Range("H125").Value = 80000
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840
